# Update the daily COVID-19 Valais figures for rows 233-254
# (2020-11-13 .. 2020-12-04).
#
# Columns B (Cumul cas positifs), H (Total hospitalisations), J (Cumul
# deces) and K (Nb nouveaux deces) are all driven by shared formulas
# already present in the sheet ( B = prev.B + C, H = G + E,
# J = prev.J + K, K = L + M ), so only the genuine input columns
# (C, E, F, G, L, M) need to be written here - Excel's own recalculation
# keeps the cumulative/derived columns in sync automatically.
#
# Columns L and M are formatted as Text ("@"), so a plain
# `Range.Value = <number>` would store the figure as text instead of a
# number (matching how Excel itself treats a Text-formatted cell).  The
# source file stores these as real numbers, so Set-NumericValue below
# briefly switches the cell to General formatting, assigns the value, and
# restores the original (Text) number format - this mirrors what a user
# re-formatting the cell before typing would produce.

function Set-NumericValue($range, $value) {
    $oldFormat = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $oldFormat
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 233 (2020-11-13): new positive cases revised up by 1 ---
$ws.Range("C233").Value = 311

# --- Row 245 (2020-11-25): new positive cases revised down by 1 ---
$ws.Range("C245").Value = 887

# --- Row 247 (2020-11-27): new cases down by 1, one fewer hospital death ---
$ws.Range("C247").Value = 784
Set-NumericValue $ws.Range("M247") 4

# --- Row 250 (2020-11-30): two more extra-hospital deaths ---
Set-NumericValue $ws.Range("M250") 5

# --- Row 251 (2020-12-01): new cases up, more hospital deaths ---
$ws.Range("C251").Value = 850
Set-NumericValue $ws.Range("L251") 3

# --- Row 252 (2020-12-02): new cases up, more deaths both in & out of hospital ---
$ws.Range("C252").Value = 643
Set-NumericValue $ws.Range("L252") 4
Set-NumericValue $ws.Range("M252") 6

# --- Row 253 (2020-12-03): new cases up, more deaths both in & out of hospital ---
$ws.Range("C253").Value = 320
Set-NumericValue $ws.Range("L253") 2
Set-NumericValue $ws.Range("M253") 2

# --- Row 254 (2020-12-04): day's figures added (was previously blank) ---
$ws.Range("C254").Value = 18
$ws.Range("E254").Value = 30
$ws.Range("F254").Value = 16
$ws.Range("G254").Value = 298

# Best-effort: restore the scroll position of the frozen pane so the view
# that was left open shows row 225 at the top (cosmetic only).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 225
$aw.ScrollColumn = 2
